$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the transfer date: the transfer now took place on Monday 15 April
# 2024 instead of Friday 5 April 2024.
$ws.Range("D10").Value = "Lunes 15 de Abril de 2024"

# Update the transfer purpose/description: memory dump report now also
# covers the disk image of the compromised corporate server.
$ws.Range("D8").Value = "Elaboración de un informe técnico sobre el volcado de memoria y la imagen del disco del servidor corporativo comprometido."

# Update the active-cell selection to match the saved view.
$ws.Range("E14").Select()
